$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text format to prevent Excel from auto-converting
# numeric-looking price strings into numbers (losing formatting like trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "60.005.26"
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("D3").Value = "2.544.74"
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "538.14"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "144.00"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").Value = "2.565.05"
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  +4.10%  "
$ws.Range("D14").Value = "2.994.96"
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").Value = "24.08"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "59.963.66"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("E17").Value = "  +4.79%  "
$ws.Range("D18").Value = "2.552.70"
$ws.Range("E18").Value = "  +2.99%  "
$ws.Range("D19").Value = "11.25"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").Value = "326.99"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "5.94"
$ws.Range("E23").Value = "  +4.13%  "
$ws.Range("D24").Value = "63.11"
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("D25").Value = "0.434"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("E26").Value = "  +4.48%  "
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  +4.16%  "
$ws.Range("D29").Value = "7.07"
$ws.Range("E29").Value = "  +4.89%  "
$ws.Range("D30").Value = "0.0₃0795"
$ws.Range("E30").Value = "  +5.15%  "
$ws.Range("D31").Value = "1.81"
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").Value = "164.99"
$ws.Range("E33").Value = "  +5.26%  "
$ws.Range("E34").Value = "  +5.82%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").Value = "18.73"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D37").Value = "4.40"
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("D39").Value = "36.99"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").Value = "5.58"
$ws.Range("E40").Value = "  -4.22%  "
$ws.Range("D41").Value = "300.33"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").Value = "0.837"
$ws.Range("E42").Value = "  +7.56%  "
$ws.Range("D43").Value = "3.72"
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "0.608"
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("D46").Value = "10.82"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "126.94"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("D48").Value = "18.94"
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("D49").Value = "0.0937"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").Value = "0.0519"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").Value = "0.0229"
$ws.Range("E51").Value = "  +1.50%  "

# Restore default (Normal) style so unchanged formatting/style indices match the original workbook.
$ws.Range("D2:E51").Style = "Normal"
